# Update bus voltage-magnitude results (res_bus/vm_pu) for the 380 kV case.
# Slack-bus voltage setpoint (column B) drops from 1.05 pu to 1.02 pu, and the
# downstream bus voltages (columns C-F, I-N) are refreshed with the re-run
# power-flow results for rows 2-25 (time steps 0-23).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$vmPuUpdates = @{
    2 = @{ "B2"=1.02; "C2"=1.040247152113656; "D2"=1.041652388748826; "E2"=1.053691823810803; "F2"=1.060856199127629; "I2"=1.036238780645264; "J2"=1.045334790692028; "K2"=1.044431297440803; "L2"=1.056437092122106; "M2"=1.063581841326341; "N2"=1.046819285996387 }
    3 = @{ "B3"=1.02; "C3"=1.041363228920941; "D3"=1.042480958939132; "E3"=1.054791697451047; "F3"=1.062020729773801; "I3"=1.036447934549956; "J3"=1.046095574662389; "K3"=1.045070492057994; "L3"=1.057349326768156; "M3"=1.064559996284495; "N3"=1.04758115036725 }
    4 = @{ "B4"=1.02; "C4"=1.04208553379342; "D4"=1.043016956478672; "E4"=1.055503870959546; "F4"=1.062774784320758; "I4"=1.036581890282143; "J4"=1.046587431556944; "K4"=1.04548330457556; "L4"=1.057939492770725; "M4"=1.065192872156668; "N4"=1.048073705755011 }
    5 = @{ "B5"=1.02; "C5"=1.042389222224798; "D5"=1.043242255373739; "E5"=1.055803384258052; "F5"=1.0630919149588; "I5"=1.036637874583791; "J5"=1.046794107494334; "K5"=1.045656661952568; "L5"=1.058187572225359; "M5"=1.065458920059114; "N5"=1.048280675195938 }
    6 = @{ "B6"=1.02; "C6"=1.042440214713822; "D6"=1.043280082019795; "E6"=1.055853680641809; "F6"=1.063145169993176; "I6"=1.036647255199347; "J6"=1.046828803392478; "K6"=1.045685758333193; "L6"=1.058229224318961; "M6"=1.065503589910685; "N6"=1.048315420366237 }
    7 = @{ "B7"=1.02; "C7"=1.042089591568606; "D7"=1.043019967069349; "E7"=1.055507872616553; "F7"=1.062779021339927; "I7"=1.036582639646643; "J7"=1.046590193566452; "K7"=1.045485621728087; "L7"=1.057942807723764; "M7"=1.065196427152987; "N7"=1.048076471686889 }
    8 = @{ "B8"=1.02; "C8"=1.040624309292877; "D8"=1.041932437314174; "E8"=1.054063431653748; "F8"=1.061249649137317; "I8"=1.03630975102764; "J8"=1.045591988492711; "K8"=1.044647479601687; "L8"=1.056745409154277; "M8"=1.063912425072475; "N8"=1.047076849047448 }
    9 = @{ "B9"=1.02; "C9"=1.03804324813765; "D9"=1.040014989681512; "E9"=1.05152182066036; "F9"=1.058558707934009; "I9"=1.035818315779816; "J9"=1.043829791138343; "K9"=1.04316452543417; "L9"=1.054634581977514; "M9"=1.061649410736712; "N9"=1.04531214917071 }
    10 = @{ "B10"=1.02; "C10"=1.036323135631786; "D10"=1.038735978299049; "E10"=1.049829868411293; "F10"=1.05676741715381; "I10"=1.035483590663965; "J10"=1.042652807980692; "K10"=1.042171829942078; "L10"=1.053226771459113; "M10"=1.060140420984568; "N10"=1.044133494561973 }
    11 = @{ "B11"=1.02; "C11"=1.035578437946352; "D11"=1.038181985282634; "E11"=1.049097811736731; "F11"=1.055992395301639; "I11"=1.03533696736881; "J11"=1.042142639692459; "K11"=1.041741018658918; "L11"=1.052617029128767; "M11"=1.059486932362811; "N11"=1.043622601776254 }
    12 = @{ "B12"=1.02; "C12"=1.035301841396609; "D12"=1.03797618174516; "E12"=1.048825978365334; "F12"=1.05570460999138; "I12"=1.035282251782993; "J12"=1.041953060943185; "K12"=1.04158085066976; "L12"=1.052390520505053; "M12"=1.059244184116304; "N12"=1.043432753803411 }
    13 = @{ "B13"=1.02; "C13"=1.035361171528286; "D13"=1.038020328471835; "E13"=1.048884283705158; "F13"=1.055766336750643; "I13"=1.035293999918106; "J13"=1.04199372983359; "K13"=1.041615213833236; "L13"=1.052439108430006; "M13"=1.059296255057682; "N13"=1.043473480448305 }
    14 = @{ "B14"=1.02; "C14"=1.035555574041884; "D14"=1.038164974025237; "E14"=1.049075340173098; "F14"=1.055968604998795; "I14"=1.035332449725898; "J14"=1.042126970678547; "K14"=1.041727782091801; "L14"=1.052598306330743; "M14"=1.059466866995051; "N14"=1.043606910510545 }
    15 = @{ "B15"=1.02; "C15"=1.035675354111162; "D15"=1.038254091499446; "E15"=1.049193067625549; "F15"=1.056093241354716; "I15"=1.035356106374692; "J15"=1.042209054210414; "K15"=1.041797119827553; "L15"=1.052696390326689; "M15"=1.059571984838753; "N15"=1.043689110610442 }
    16 = @{ "B16"=1.02; "C16"=1.036372561252631; "D16"=1.038772741359874; "E16"=1.049878464561406; "F16"=1.05681886578384; "I16"=1.035493286066088; "J16"=1.042686655037646; "K16"=1.042200401071293; "L16"=1.053267234825727; "M16"=1.060183789012014; "N16"=1.044167389685631 }
    17 = @{ "B17"=1.02; "C17"=1.03680993289017; "D17"=1.039098030545074; "E17"=1.050308548128826; "F17"=1.05727419603843; "I17"=1.035578884112446; "J17"=1.042986100104684; "K17"=1.042453109601185; "L17"=1.053625269713752; "M17"=1.060567534173372; "N17"=1.044467259999011 }
    18 = @{ "B18"=1.02; "C18"=1.03706505624544; "D18"=1.039287749560546; "E18"=1.050559463612045; "F18"=1.057539842229576; "I18"=1.035628649415647; "J18"=1.043160710541921; "K18"=1.04260041686197; "L18"=1.053834090860709; "M18"=1.060791357962216; "N18"=1.044642118403096 }
    19 = @{ "B19"=1.02; "C19"=1.037152048763672; "D19"=1.039352436013464; "E19"=1.050645028655343; "F19"=1.057630430871962; "I19"=1.035645590519973; "J19"=1.04322023958309; "K19"=1.042650629001376; "L19"=1.053905291022622; "M19"=1.060867674703562; "N19"=1.044701731982332 }
    20 = @{ "B20"=1.02; "C20"=1.036763005831476; "D20"=1.039063131816872; "E20"=1.050262398545586; "F20"=1.057225337234982; "I20"=1.035569717069438; "J20"=1.042953977740469; "K20"=1.042426006023901; "L20"=1.053586857466495; "M20"=1.060526362808575; "N20"=1.044435092017354 }
    21 = @{ "B21"=1.02; "C21"=1.035498326877126; "D21"=1.038122380228365; "E21"=1.049019076461135; "F21"=1.055909039498156; "I21"=1.03532113420442; "J21"=1.042087736769197; "K21"=1.041694637562798; "L21"=1.052551427160389; "M21"=1.059416626408218; "N21"=1.043567620884544 }
    22 = @{ "B22"=1.02; "C22"=1.034703272757409; "D22"=1.03753074275643; "E22"=1.048237842456731; "F22"=1.055081964635747; "I22"=1.035163375554901; "J22"=1.041542636353164; "K22"=1.04123395532025; "L22"=1.051900276721745; "M22"=1.058718812854276; "N22"=1.043021746363396 }
    23 = @{ "B23"=1.02; "C23"=1.035124736888478; "D23"=1.037844395131787; "E23"=1.048651942932874; "F23"=1.055520362296188; "I23"=1.035247145238761; "J23"=1.041831648160201; "K23"=1.041478251535379; "L23"=1.052245476778867; "M23"=1.05908874461822; "N23"=1.043311168600349 }
    24 = @{ "B24"=1.02; "C24"=1.036784210112298; "D24"=1.039078901101313; "E24"=1.050283251385533; "F24"=1.057247414240238; "I24"=1.035573859763893; "J24"=1.042968492610947; "K24"=1.042438253252106; "L24"=1.053604214352312; "M24"=1.060544966401436; "N24"=1.044449627500613 }
    25 = @{ "B25"=1.02; "C25"=1.038710406081347; "D25"=1.040510822535422; "E25"=1.052178453968251; "F25"=1.059253908315233; "I25"=1.035946615681246; "J25"=1.044285745650272; "K25"=1.043548620268348; "L25"=1.055180384670186; "M25"=1.062234508083593; "N25"=1.045768751190344 }
}

foreach ($row in $vmPuUpdates.Keys) {
    $rowUpdates = $vmPuUpdates[$row]
    foreach ($cellRef in $rowUpdates.Keys) {
        $ws.Range($cellRef).Value = $rowUpdates[$cellRef]
    }
}
